# Fruta / hortaliza, semanal
# The data rows (2-33) for "Hortaliza, Vega Modelo de Temuco - Ramas de apio"
# get reshuffled: the Fecha (D), Volumen (J), Precio minimo (K), Precio maximo
# (L), Precio promedio ponderado (M) and Precio $/Kg (P) values move between
# rows according to the mapping below (target row -> source row), while every
# other column (market/category/unit descriptors, etc.) stays put because it
# is identical across rows anyway.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$mapping = @{
    2  = 14
    3  = 10
    4  = 25
    5  = 15
    6  = 32
    7  = 12
    8  = 22
    9  = 3
    10 = 16
    11 = 8
    12 = 33
    13 = 27
    14 = 24
    15 = 17
    16 = 6
    17 = 26
    18 = 13
    19 = 29
    20 = 2
    21 = 7
    22 = 5
    23 = 23
    24 = 9
    25 = 20
    26 = 19
    27 = 21
    28 = 4
    29 = 30
    30 = 31
    31 = 11
    32 = 28
    33 = 18
}

$cols = @("D", "J", "K", "L", "M", "P")

# Snapshot every source row's values up front so writes to one row never
# clobber data another row still needs to read.
$snapshot = @{}
foreach ($row in $mapping.Keys) {
    $values = @{}
    foreach ($col in $cols) {
        $values[$col] = $ws.Range("$col$row").Value2
    }
    $snapshot[$row] = $values
}

foreach ($row in $mapping.Keys) {
    $srcRow = $mapping[$row]
    $srcValues = $snapshot[$srcRow]
    foreach ($col in $cols) {
        $ws.Range("$col$row").Value2 = $srcValues[$col]
    }
}
